# Generate Report for Handoff
#
# This reproduces a re-run of the localization-status report generation:
#  - The "Latest HO Xliff Generate Date" (Overview sheet) and the de-de
#    "Latest Handoff Datetime" column both pick up the new de-de handoff
#    timestamp (2016-08-20 14:17:51) for the six rows that were still
#    "Ready for handoff".
#  - The zh-cn "Latest Handoff Datetime" column picks up its own refreshed
#    handoff timestamp (2016-08-20 14:17:46) for the same rows.
#  - The zh-cn and de-de "Priority" column gets populated with "ht" for
#    those same rows (previously blank).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 11, 12, 13, 14)

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $overview.Range("G$r").Value = "2016-08-20 14:17:51"

    # zh-cn sheet: column H = "Latest Handoff Datetime"
    $zhcn.Range("H$r").Value = "2016-08-20 14:17:46"

    # zh-cn sheet: column E = "Priority"
    $zhcn.Range("E$r").Value = "ht"

    # de-de sheet: column H = "Latest Handoff Datetime" (mirrors Overview's date)
    $dede.Range("H$r").Value = "2016-08-20 14:17:51"

    # de-de sheet: column E = "Priority"
    $dede.Range("E$r").Value = "ht"
}
